# "setting spring 2018 to current"
#
# The sheet tracks semesters with a rolling "status" (current/past/future)
# and a "priority" weight. Spring_2018 (row 2) was "future"; it is now
# "current". The row that used to be "current" (Fall_2017, row 3) becomes
# "past" and gets a lastmod date. The priority values cascade down one row,
# and the now-unused "future" shared string is no longer referenced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Spring_2018: future -> current, priority 0.6
$ws.Range("B2").Value = "current"
$ws.Range("D2").Value = 0.6

# Row 3 - Fall_2017: current -> past, gets a lastmod date, priority drops to 0.3
$ws.Range("B3").Value = "past"
$ws.Range("C3").Value = 43079
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats - match the date format used below
$ws.Range("D3").Value = 0.3

# Row 4 - Spring_2017: priority drops to 0.1
$ws.Range("D4").Value = 0.1

# Row 5 - Fall_2016: priority drops to 0.05
$ws.Range("D5").Value = 0.05

# Update the active selection to match
$ws.Range("G5").Select()
